# Generate Report for Handoff
#
# The 4c675ee6-6f80-4b42-b109-b0a342d14def.md file has been fully handed
# off and is no longer pending, so its row is removed from every sheet.
# The still-pending 3a249cac-... file moves from "Handed back: in sync
# with en-US" to "Ready for handoff", with refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

# --- Remove the row for 4c675ee6-6f80-4b42-b109-b0a342d14def (row 3) ---
# from all three sheets. This shifts the dimension/rows/shared strings
# automatically.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Rows("3:3").Delete()

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Rows("3:3").Delete()

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Rows("3:3").Delete()

# --- Update status text + handoff timestamps for the remaining row ---

# Overview sheet: B2/C2 = Status columns ("zh-cn"/"de-de"), D2 = Latest
# Handoff Date.
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-43-20 16:43:55"

# zh-cn sheet: C2 = Status, E2 = Latest Handoff Datetime.
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-20 16:43:52"

# de-de sheet: C2 = Status, E2 = Latest Handoff Datetime.
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-20 16:43:55"
